$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137, shifting existing rows 137-221 down to 138-222
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new data entry
$ws.Cells.Item(137, 1).Value = 8
$ws.Cells.Item(137, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(137, 3).Value = "Coquimbo"
$ws.Cells.Item(137, 4).Value = 45029
$ws.Cells.Item(137, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(137, 5).Value = 4
$ws.Cells.Item(137, 6).Value = 100112040
$ws.Cells.Item(137, 7).Value = "Cilantro"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 2200
$ws.Cells.Item(137, 11).Value = 1800
$ws.Cells.Item(137, 12).Value = 2000
$ws.Cells.Item(137, 13).Value = 1900
$ws.Cells.Item(137, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(137, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(137, 16).Value = 1267
$ws.Cells.Item(137, 17).Value = 1.5
$ws.Cells.Item(137, 18).Value = "Hortaliza"
